$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J (DATE_TYPE_CODE) updated - force text so leading zero is kept ("004", not 4)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "004"

# Report date (REPORT_DATE) updated
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# Updated numeric metrics for row 2
$ws.Range("O2").Value = 27700941.21
$ws.Range("P2").Value = 241.1534221946
$ws.Range("Q2").Value = 206490364.59
$ws.Range("R2").Value = 1797.6233259943
$ws.Range("S2").Value = 44911247.26
$ws.Range("T2").Value = 390.9795298893
$ws.Range("U2").Value = -9712795.08
$ws.Range("V2").Value = -84.55574685569999
$ws.Range("W2").Value = 298816.02
$ws.Range("X2").Value = 2.601373913
$ws.Range("Y2").Value = 3024071.28
$ws.Range("Z2").Value = 26.3263667687
$ws.Range("AA2").Value = -29475000
$ws.Range("AB2").Value = -256.5976753389
$ws.Range("AC2").Value = -11486853.87
$ws.Range("AD2").Value = 46.5798962312
